$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-11 and add new rows 12-21 with fresh fake data
$ws.Range("A2").Value = "sbullock@example.org"
$ws.Range("B2").Value = 45528
$ws.Range("C2").Value = 236.5
$ws.Range("D2").Value = "Water"
$ws.Range("E2").Value = 14
$ws.Range("F2").Value = "retail"

$ws.Range("A3").Value = "danaacevedo@example.org"
$ws.Range("B3").Value = 45687
$ws.Range("C3").Value = 209.68
$ws.Range("D3").Value = "When"
$ws.Range("E3").Value = 33
$ws.Range("F3").Value = "online"

$ws.Range("A4").Value = "robertskelly@example.net"
$ws.Range("B4").Value = 45638
$ws.Range("C4").Value = 68.23
$ws.Range("D4").Value = "Why"
$ws.Range("E4").Value = 39
$ws.Range("F4").Value = "b2b"

$ws.Range("A5").Value = "devin10@example.net"
$ws.Range("B5").Value = 45441
$ws.Range("C5").Value = 915.1900000000001
$ws.Range("D5").Value = "Tough"
$ws.Range("E5").Value = 54
$ws.Range("F5").Value = "b2b"

$ws.Range("A6").Value = "tinataylor@example.com"
$ws.Range("B6").Value = 45373
$ws.Range("C6").Value = 335.82
$ws.Range("D6").Value = "Staff"
$ws.Range("E6").Value = 86
$ws.Range("F6").Value = "b2b"

$ws.Range("A7").Value = "boneal@example.org"
$ws.Range("B7").Value = 45643
$ws.Range("C7").Value = 787.5
$ws.Range("D7").Value = "Money"
$ws.Range("E7").Value = 72
$ws.Range("F7").Value = "online"

$ws.Range("A8").Value = "xsmith@example.org"
$ws.Range("B8").Value = 45421
$ws.Range("C8").Value = 132.95
$ws.Range("D8").Value = "Pressure"
$ws.Range("E8").Value = 82
$ws.Range("F8").Value = "b2b"

$ws.Range("A9").Value = "davidharding@example.com"
$ws.Range("B9").Value = 45509
$ws.Range("C9").Value = 326.88
$ws.Range("D9").Value = "News"
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = "retail"

$ws.Range("A10").Value = "cjohnson@example.net"
$ws.Range("B10").Value = 45584
$ws.Range("C10").Value = 156.27
$ws.Range("D10").Value = "Daughter"
$ws.Range("E10").Value = 74
$ws.Range("F10").Value = "retail"

$ws.Range("A11").Value = "tlove@example.com"
$ws.Range("B11").Value = 45590
$ws.Range("C11").Value = 181.99
$ws.Range("D11").Value = "Fact"
$ws.Range("E11").Value = 48
$ws.Range("F11").Value = "online"

$ws.Range("A12").Value = "williamsmichelle@example.net"
$ws.Range("B12").Value = 45379
$ws.Range("C12").Value = 176.09
$ws.Range("D12").Value = "Before"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = "online"

$ws.Range("A13").Value = "rthomas@example.net"
$ws.Range("B13").Value = 45651
$ws.Range("C13").Value = 488.08
$ws.Range("D13").Value = "Technology"
$ws.Range("E13").Value = 14
$ws.Range("F13").Value = "b2b"

$ws.Range("A14").Value = "hawkinsjohnny@example.net"
$ws.Range("B14").Value = 45695
$ws.Range("C14").Value = 922.24
$ws.Range("D14").Value = "Add"
$ws.Range("E14").Value = 93
$ws.Range("F14").Value = "b2b"

$ws.Range("A15").Value = "john98@example.net"
$ws.Range("B15").Value = 45616
$ws.Range("C15").Value = 975.3099999999999
$ws.Range("D15").Value = "Often"
$ws.Range("E15").Value = 94
$ws.Range("F15").Value = "retail"

$ws.Range("A16").Value = "mayerdiana@example.org"
$ws.Range("B16").Value = 45575
$ws.Range("C16").Value = 310.58
$ws.Range("D16").Value = "Interesting"
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = "online"

$ws.Range("A17").Value = "julie52@example.org"
$ws.Range("B17").Value = 45558
$ws.Range("C17").Value = 219.21
$ws.Range("D17").Value = "Full"
$ws.Range("E17").Value = 35
$ws.Range("F17").Value = "retail"

$ws.Range("A18").Value = "davidwalker@example.com"
$ws.Range("B18").Value = 45573
$ws.Range("C18").Value = 404.18
$ws.Range("D18").Value = "Partner"
$ws.Range("E18").Value = 29
$ws.Range("F18").Value = "online"

$ws.Range("A19").Value = "tdiaz@example.org"
$ws.Range("B19").Value = 45570
$ws.Range("C19").Value = 452.66
$ws.Range("D19").Value = "Team"
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = "b2b"

$ws.Range("A20").Value = "alyssa08@example.org"
$ws.Range("B20").Value = 45501
$ws.Range("C20").Value = 411.06
$ws.Range("D20").Value = "Early"
$ws.Range("E20").Value = 66
$ws.Range("F20").Value = "retail"

$ws.Range("A21").Value = "fishersamuel@example.org"
$ws.Range("B21").Value = 45480
$ws.Range("C21").Value = 263.35
$ws.Range("D21").Value = "Responsibility"
$ws.Range("E21").Value = 31
$ws.Range("F21").Value = "b2b"

# Apply the date number format used by column B to the full range (matches existing style)
$ws.Range("B2:B21").NumberFormat = "YYYY-MM-DD"